$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------- 1. Reset sheet ----------
$ws.Range("A5:A10").UnMerge()
$ws.Range("A11:A14").UnMerge()
$ws.Range("A1:C1").UnMerge()
$ws.Cells.Clear()

# ---------- 2. Style definitions (index -> properties) ----------
function Set-CellStyle($rng, [int]$idx) {
    switch ($idx) {
        0 {
            $rng.Font.Name = "Calibri"
            $rng.Font.Size = 12
            $rng.Font.Color = 0
            $rng.Interior.Pattern = -4142
            $rng.HorizontalAlignment = 1
            $rng.VerticalAlignment = -4107
            $rng.WrapText = $false
        }
        1 {
            $rng.Font.Name = "Futura Medium"
            $rng.Font.Size = 18
            $rng.Font.Color = 0
            $rng.Interior.Pattern = -4142
            $rng.HorizontalAlignment = 1
            $rng.VerticalAlignment = -4107
            $rng.WrapText = $false
        }
        2 {
            $rng.Font.Name = "Futura Medium"
            $rng.Font.Size = 18
            $rng.Font.Color = 65535
            $rng.Interior.Pattern = 1
            $rng.Interior.Color = 0
            $rng.HorizontalAlignment = 1
            $rng.VerticalAlignment = -4107
            $rng.WrapText = $false
        }
        3 {
            $rng.Font.Name = "Futura Medium"
            $rng.Font.Size = 18
            $rng.Font.Color = 0
            $rng.Interior.Pattern = 1
            $rng.Interior.Color = 65535
            $rng.HorizontalAlignment = 1
            $rng.VerticalAlignment = -4107
            $rng.WrapText = $false
        }
        4 {
            $rng.Font.Name = "Futura Medium"
            $rng.Font.Size = 12
            $rng.Font.Color = 0
            $rng.Interior.Pattern = -4142
            $rng.HorizontalAlignment = 1
            $rng.VerticalAlignment = -4107
            $rng.WrapText = $false
        }
        5 {
            $rng.Font.Name = "Futura Medium"
            $rng.Font.Size = 12
            $rng.Font.Color = 0
            $rng.Interior.Pattern = -4142
            $rng.HorizontalAlignment = 1
            $rng.VerticalAlignment = -4107
            $rng.WrapText = $true
        }
        6 {
            $rng.Font.Name = "Futura Medium"
            $rng.Font.Size = 18
            $rng.Font.Color = 0
            $rng.Interior.Pattern = -4142
            $rng.HorizontalAlignment = -4131
            $rng.VerticalAlignment = -4108
            $rng.WrapText = $true
        }
        7 {
            $rng.Font.Name = "Futura Medium"
            $rng.Font.Size = 18
            $rng.Font.Color = 0
            $rng.Interior.Pattern = -4142
            $rng.HorizontalAlignment = -4131
            $rng.VerticalAlignment = -4108
            $rng.WrapText = $false
        }
        8 {
            $rng.Font.Name = "Futura Medium"
            $rng.Font.Size = 18
            $rng.Font.Color = 0
            $rng.Interior.Pattern = -4142
            $rng.HorizontalAlignment = -4131
            $rng.VerticalAlignment = -4108
            $rng.WrapText = $true
        }
        9 {
            $rng.Font.Name = "Futura Medium"
            $rng.Font.Size = 16
            $rng.Font.Color = 65535
            $rng.Interior.Pattern = 1
            $rng.Interior.Color = 0
            $rng.HorizontalAlignment = -4108
            $rng.VerticalAlignment = -4107
            $rng.WrapText = $true
        }
    }
}

# ---------- 3. Cell data: Ref|StyleIndex|Text ("@@EMPTY@@" = no value) ----------
$cellData = @'
A1|9|PROCESO PORTAFOLIO desafío fina móduo 1 Bootcamp DESAFIO LATAM
B1|9|@@EMPTY@@
C1|9|@@EMPTY@@
D1|1|@@EMPTY@@
E1|1|@@EMPTY@@
F1|1|@@EMPTY@@
G1|1|@@EMPTY@@
H1|1|@@EMPTY@@
A2|2|Tareas 
B2|2|Subtarea
C2|2|estado 
D2|1|@@EMPTY@@
E2|1|@@EMPTY@@
F2|1|@@EMPTY@@
G2|1|@@EMPTY@@
H2|1|@@EMPTY@@
A3|1|Definición template 
B3|1|@@EMPTY@@
C3|3|ok
D3|1|@@EMPTY@@
E3|1|@@EMPTY@@
F3|1|@@EMPTY@@
G3|1|@@EMPTY@@
H3|1|@@EMPTY@@
A4|1|Producción de textos 
B4|1|@@EMPTY@@
C4|3|ok
D4|1|@@EMPTY@@
E4|1|@@EMPTY@@
F4|1|@@EMPTY@@
G4|1|@@EMPTY@@
H4|1|@@EMPTY@@
A5|7|Edición de imágenes e incorporación al sitio
B5|1|Desafio cv
C5|3|ok
D5|1|@@EMPTY@@
E5|4|Github 
F5|1|@@EMPTY@@
G5|1|@@EMPTY@@
H5|1|@@EMPTY@@
A6|7|@@EMPTY@@
B6|1|Desafio landing page Caps
C6|3|ok
D6|1|@@EMPTY@@
E6|4|paso a paso de examen final del primer módulo
F6|1|@@EMPTY@@
G6|1|@@EMPTY@@
H6|1|@@EMPTY@@
A7|7|@@EMPTY@@
B7|1|Desafio Iguana page
C7|3|ok
D7|1|@@EMPTY@@
E7|4|// esta el parte de creación de repos x desafio
F7|1|@@EMPTY@@
G7|1|@@EMPTY@@
H7|1|@@EMPTY@@
A8|7|@@EMPTY@@
B8|1|Desafio Surfon_Cupon
C8|3|ok
D8|1|@@EMPTY@@
E8|4|@@EMPTY@@
F8|1|@@EMPTY@@
G8|1|@@EMPTY@@
H8|1|@@EMPTY@@
A9|7|@@EMPTY@@
B9|1|Bonus track Efdun
C9|3|ok
D9|1|@@EMPTY@@
E9|4|@@EMPTY@@
F9|1|@@EMPTY@@
G9|1|@@EMPTY@@
H9|1|@@EMPTY@@
A10|7|@@EMPTY@@
B10|1|Bonus track El Fuego
C10|3|ok
D10|1|@@EMPTY@@
E10|4|@@EMPTY@@
F10|1|@@EMPTY@@
G10|1|@@EMPTY@@
H10|1|@@EMPTY@@
A11|8|Repositorios (creación repos remotos, carga con git, obtención línks)
B11|1|Desafio cv
C11|3|ok
D11|4|@@EMPTY@@
E11|4|@@EMPTY@@
F11|1|@@EMPTY@@
G11|1|@@EMPTY@@
H11|1|@@EMPTY@@
A12|8|@@EMPTY@@
B12|1|Desafio landing page Caps
C12|3|ok
E12|4|@@EMPTY@@
A13|8|@@EMPTY@@
B13|1|Desafio Iguana page
C13|3|ok
E13|4|1-creo repositorio en github con el nombre correspondiente
A14|8|@@EMPTY@@
B14|1|Desafio Surfon_Cupon
C14|3|ok
E14|5|2-en mi local, voy a visual studio code y abro la carpeta que quiero subir al archivo recién creado
A15|1|Forkeos a dos repos 
B15|4|forkeo repo en user de github (seccion repositorios)
C15|3|@@EMPTY@@
E15|5|3-git  init > add . >  status > commit -m "mesage" > establezco remote desde la terminal del proyecto >sigo proceso que github me indica en terminal > git push
A16|6|@@EMPTY@@
B16|4|abro carpeta de proyecto en visual studio code
C16|3|@@EMPTY@@
E16|4|4-check que haya estabecido conección y que se subieron archivos de local a remoto
A17|6|@@EMPTY@@
B17|5|en terminal de carpeta padre escribo git clone + enlace ssh que obtengo del repo forkeado
C17|3|@@EMPTY@@
A18|6|@@EMPTY@@
B18|5|en terminal de carpeta padre escribo git clone + enlace ssh que obtengo del repo forkeado
C18|3|@@EMPTY@@
A19|6|@@EMPTY@@
B19|5|en terminal de carpeta padre escribo git clone + enlace ssh que obtengo del repo forkeado
C19|3|@@EMPTY@@
B20|5|en terminal  cd a carpeta del proyecto ya clonado
C20|1|@@EMPTY@@
A21|1|@@EMPTY@@
B21|5|abro el archivo a modificar (html / css / img) hago cambios en index o css
C21|1|@@EMPTY@@
A22|1|@@EMPTY@@
B22|5|git add . > git commit -m "" > git log > git push origin master / main (según lo que indique la consola)
C22|1|@@EMPTY@@
A23|1|@@EMPTY@@
B23|5|chekeo en repo forkeado
C23|1|@@EMPTY@@
A24|1|@@EMPTY@@
B24|1|@@EMPTY@@
C24|1|@@EMPTY@@
A25|1|Check de reuqrimientos examen
B25|1|@@EMPTY@@
C25|1|@@EMPTY@@
A26|1|@@EMPTY@@
B26|1|@@EMPTY@@
C26|1|@@EMPTY@@
D26|4|@@EMPTY@@
A27|1|detalle 
B27|1|direccionamiento links nav bar / cómo se hace anchor a partes del portafolio, sapear trabajo de roberto de guiro
C27|1|@@EMPTY@@
D27|4|@@EMPTY@@
A28|1|@@EMPTY@@
B28|1|@@EMPTY@@
C28|1|@@EMPTY@@
D28|4|@@EMPTY@@
A29|4|@@EMPTY@@
B29|4|@@EMPTY@@
C29|4|@@EMPTY@@
D29|4|@@EMPTY@@
A30|4|ideas sitio enchuado
B30|4|@@EMPTY@@
C30|4|@@EMPTY@@
D30|4|@@EMPTY@@
A31|4|frase > skills > link a inkdin
B31|4|@@EMPTY@@
C31|4|@@EMPTY@@
D31|4|@@EMPTY@@

'@

$rows = $cellData -split "`n"
foreach ($line in $rows) {
    $line = $line.Trim("`r")
    if ($line.Length -eq 0) { continue }
    $parts = $line -split "\|", 3
    $ref = $parts[0]
    $styleIdx = [int]$parts[1]
    $text = $parts[2]
    $rng = $ws.Range($ref)
    if ($text -ne "@@EMPTY@@") {
        $rng.Value = $text
    }
    Set-CellStyle $rng $styleIdx
}

# ---------- 4. Row heights ----------
$heightData = @'
1=29
2=26
3=26
4=26
5=26
6=26
7=26
8=26
9=26
10=26
11=47
12=26
13=26
14=54
15=36
16=26
17=37
18=37
19=37
20=26
21=37
22=37
23=26
24=37
25=26
26=26
27=26
28=26

'@
$hlines = $heightData -split "`n"
foreach ($line in $hlines) {
    $line = $line.Trim("`r")
    if ($line.Length -eq 0) { continue }
    $hp = $line -split "="
    $rowNum = [int]$hp[0]
    $rowH = [double]$hp[1]
    $ws.Rows.Item($rowNum).RowHeight = $rowH
}

# ---------- 5. Merges ----------
$ws.Range("A5:A10").Merge()
$ws.Range("A11:A14").Merge()
$ws.Range("A1:C1").Merge()

# ---------- 6. View / selection ----------
[void]$ws.Range("B24").Select()
